# Re-order the worker/period rows in the "Estado de Cuenta" table.
#
# Before:
#   Row16: CC 1235038025 RICARDO MARIO JIMENEZ RESTREPO  2409  52000
#   Row17: CC 13541643   PEDRO ALONSO HERNANDEZ ROMERO    2409  52000
#   Row18: CC 1235038025 RICARDO MARIO JIMENEZ RESTREPO  2410  15600
#   Row19: CC 13541643   PEDRO ALONSO HERNANDEZ ROMERO    2410  15600
#
# After (grouped by worker instead of by period):
#   Row16: CC 1235038025 RICARDO MARIO JIMENEZ RESTREPO  2410  15600
#   Row17: CC 1235038025 RICARDO MARIO JIMENEZ RESTREPO  2409  52000
#   Row18: CC 13541643   PEDRO ALONSO HERNANDEZ ROMERO    2410  15600
#   Row19: CC 13541643   PEDRO ALONSO HERNANDEZ ROMERO    2409  52000

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2410"
$ws.Range("F16").Value = 15600

$ws.Range("C17").Value = "1235038025"
$ws.Range("D17").Value = "RICARDO MARIO JIMENEZ RESTREPO"

$ws.Range("C18").Value = "13541643"
$ws.Range("D18").Value = "PEDRO ALONSO HERNANDEZ ROMERO"

$ws.Range("E19").Value = "2409"
$ws.Range("F19").Value = 52000
